$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.253.41'
$ws.Range('E2').Value = '  -2.37%  '
$ws.Range('D3').Value = '1.559.44'
$ws.Range('E3').Value = '  -3.68%  '
$ws.Range('E4').Value = '  -0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '206.11'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.25%  '
$ws.Range('E6').Value = '  -0.24%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.477'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -4.42%  '
$ws.Range('E8').Value = '  -0.57%  '
$ws.Range('E9').Value = '  -2.61%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '17.74'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.29%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0782'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.92%  '
$ws.Range('D12').Value = '1.776.40'
$ws.Range('E12').Value = '  -3.64%  '
$ws.Range('D13').Value = '1.584.31'
$ws.Range('E13').Value = '  -2.14%  '
$ws.Range('E14').Value = '  -3.28%  '
$ws.Range('E15').Value = '  -2.95%  '
$ws.Range('D16').Value = '25.263.56'
$ws.Range('E16').Value = '  -2.32%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').Value = '0.0₃0712'
$ws.Range('E17').Value = '  -3.44%  '
$ws.Range('B18').Value = 'Litecoin'
$ws.Range('C18').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '59.18'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.54%  '
$ws.Range('E19').Value = '  -0.26%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '186.39'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.50%  '
$ws.Range('E21').Value = '  -2.83%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.25'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.85'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.58%  '
$ws.Range('E24').Value = '  -2.49%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '141.10'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.50%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.23%  '
$ws.Range('E27').Value = '  -2.63%  '
$ws.Range('E28').Value = '  -1.61%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.37'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.49%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.14'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -7.13%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0464'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.58%  '
$ws.Range('E32').Value = '  -2.11%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.98'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.84%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.48'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.31'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.18%  '
$ws.Range('D36').Value = '1.084.31'
$ws.Range('E36').Value = '  -3.33%  '
$ws.Range('B37').Value = 'PaxDollar'
$ws.Range('C37').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.23%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.33'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.06%  '
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.494'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.96%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0148'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.12%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.769'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -7.70%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.796'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +6.42%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '93.23'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.99%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.08'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.83%  '
$ws.Range('D45').Value = '1.692.25'
$ws.Range('E45').Value = '  -3.53%  '
$ws.Range('E46').Value = '  -0.84%  '
$ws.Range('E47').Value = '  -1.51%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '52.46'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.78%  '
$ws.Range('E49').Value = '  -3.13%  '
$ws.Range('E50').Value = '  -0.20%  '
$ws.Range('E51').Value = '  -2.11%  '
